# "update logo, doc, and fix"
#
# docs/_static/logo.pptx gains a 4th slide: a new "big logo" slide (the same
# layout/shapes already used on slides 1 and 2 of this deck) is inserted
# right before the existing small-logo slide, which is otherwise left
# untouched and simply slides down to the last position.
#
# Target p:sldIdLst goes from
#     257, 258, 259(old small-logo slide)
# to
#     257, 258, 260(new big-logo slide), 259(old small-logo slide, unchanged)

$p = $ppt.ActivePresentation

# Slide 1 already has the exact shape tree we need for the new slide
# (empty styled placeholder box, big "SAIUnit" wordmark, "Unit-aware
# Computations for Scientific AI" tagline, underline connector). Duplicate
# it and move the duplicate to slide position 3, right before the old
# third slide (which shifts down to position 4 automatically).
$srcSlide = $p.Slides.Item(1)
$dupRange = $srcSlide.Duplicate()
$newSlide = $dupRange.Item(1)
$newSlide.MoveTo(3)

# Shape 2 ("文本框 3") holds the big wordmark run - fix the casing to SAIUnit.
$wordmark = $newSlide.Shapes.Item(2)
$wordmark.TextFrame.TextRange.Runs(1).Text = "SAIUnit"

# Shape 3 ("文本框 4") holds the tagline - lower-case/singular "computation".
$tagline = $newSlide.Shapes.Item(3)
$tagline.TextFrame.TextRange.Runs(2).Text = "-aware computation for "
